$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the current column B (ASIN), shifting
# ASIN..is_holiday_week one column to the right (B->C ... I->J).
$ws.Columns.Item(2).Insert()

# New column header
$ws.Range("B1").Value = "Week_Start_Date"

# Week labels (A) lose their leading zero, and the new Week_Start_Date
# column (B) gets the corresponding week-start date, stored as plain text.
$weeks = @(
    @{ Row = 2;  Label = "W1";  Date = "2025-01-05" },
    @{ Row = 3;  Label = "W2";  Date = "2025-01-12" },
    @{ Row = 4;  Label = "W3";  Date = "2025-01-19" },
    @{ Row = 5;  Label = "W4";  Date = "2025-01-26" },
    @{ Row = 6;  Label = "W5";  Date = "2025-02-02" },
    @{ Row = 7;  Label = "W6";  Date = "2025-02-09" },
    @{ Row = 8;  Label = "W7";  Date = "2025-02-16" },
    @{ Row = 9;  Label = "W8";  Date = "2025-02-23" },
    @{ Row = 10; Label = "W9";  Date = "2025-03-02" },
    @{ Row = 11; Label = "W10"; Date = "2025-03-09" },
    @{ Row = 12; Label = "W11"; Date = "2025-03-16" },
    @{ Row = 13; Label = "W12"; Date = "2025-03-23" },
    @{ Row = 14; Label = "W13"; Date = "2025-03-30" },
    @{ Row = 15; Label = "W14"; Date = "2025-04-06" },
    @{ Row = 16; Label = "W15"; Date = "2025-04-13" },
    @{ Row = 17; Label = "W16"; Date = "2025-04-20" }
)

foreach ($w in $weeks) {
    $row = $w.Row

    $ws.Cells.Item($row, 1).Value = $w.Label

    # Force text storage so the date string isn't auto-converted to a
    # date serial number.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $w.Date

    # is_holiday_week (now column J) becomes a boolean column.
    $ws.Cells.Item($row, 10).Value = $false
}

$wb.Save()
